# Realestate Update resale numbers 2023-07-01 12:39
# Appends the new daily snapshot row (row 97) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$newRow = 97

# Columns A:D hold text (date/time/weekday/week are stored as plain text,
# matching the rest of the sheet) - force text format so Excel doesn't
# auto-convert the date/number-looking strings.
$ws.Range("A" + $newRow + ":D" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2023-07-01"
$ws.Range("B" + $newRow).Value = "12:30:12"
$ws.Range("C" + $newRow).Value = "Saturday"
$ws.Range("D" + $newRow).Value = "26"

# Columns E:T hold the numeric city resale figures.
$values = @{
    "E" = 123492
    "F" = 134817
    "G" = 161159
    "H" = 131507
    "I" = 175398
    "J" = 112884
    "K" = 204741
    "L" = 221999
    "M" = 173990
    "N" = 103257
    "O" = 38715
    "P" = 32743
    "Q" = 51887
    "R" = -1
    "S" = 35931
    "T" = -1
}

foreach ($col in $values.Keys) {
    $ws.Range($col + $newRow).Value = $values[$col]
}
